$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 131085086
$ws.Range("B33").Value = 57884
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 100109
$ws.Range("F33").Value = "Tretåig hackspett"
$ws.Range("G33").Value = "Picoides tridactylus"
$ws.Range("H33").Value = "(Linnaeus, 1758)"

# I33 is blank in the source data (present, but empty) - materialize the
# cell without leaving a stray value/style behind.
$ws.Range("I33").NumberFormat = "General"
$ws.Range("I33").Style = "Normal"

$ws.Range("M33").Value = "färska spår"
$ws.Range("P33").Value = "Sör-Tågsjöberget, Sör-Tågsjöberget, Ång"
$ws.Range("Q33").Value = 585166
$ws.Range("R33").Value = 7060188
$ws.Range("S33").Value = 15
$ws.Range("T33").Value = "Västernorrland"
$ws.Range("U33").Value = "Sollefteå"
$ws.Range("V33").Value = "Ångermanland"
$ws.Range("W33").Value = "Junsele"

$ws.Range("Y33").NumberFormat = "@"
$ws.Range("Y33").Value = "2026-02-09"
$ws.Range("Y33").Style = "Normal"
$ws.Range("AA33").NumberFormat = "@"
$ws.Range("AA33").Value = "2026-02-09"
$ws.Range("AA33").Style = "Normal"

$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AG33").Value = $false

# AT33 is blank in the source data, same as I33 above.
$ws.Range("AT33").NumberFormat = "General"
$ws.Range("AT33").Style = "Normal"

$ws.Range("AW33").Value = "Daniel Rutschman"
$ws.Range("AX33").Value = "Daniel Rutschman"

# AY33 is blank in the source data, same as I33 above.
$ws.Range("AY33").NumberFormat = "General"
$ws.Range("AY33").Style = "Normal"
